$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Active sheet / selection changes.
#    The active tab moves from "dcin5_log2_expression" (0-based
#    tab index 3) to "optimization_parameters" (0-based tab index 6),
#    and the selection on optimization_parameters becomes C1:I10.
# ---------------------------------------------------------------
$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Activate()
$wsOpt.Range("C1:I10").Select()

# ---------------------------------------------------------------
# 2. optimization_parameters sheet restructuring:
#    - Row 1: drop the superfluous "value" cells in C1:F1, keeping
#      only the A1/B1 header pair.
#    - The "Model"/"Sigmoid" row becomes "production_function"/"Sigmoid".
#    - A brand new "L_curve" row (value 1) is inserted right after it.
#    - The old "Deletion" row (0 / 3) is deleted entirely, so every
#      row that used to follow it shifts up by one.
# ---------------------------------------------------------------
$wsOpt.Range("C1:F1").ClearContents()

$wsOpt.Range("A8").Value = "production_function"

$wsOpt.Rows.Item(9).Insert()
$wsOpt.Range("A9").Value = "L_curve"
$wsOpt.Range("B9").Value = 1
$wsOpt.Range("B9").NumberFormat = "0.00E+00"

$wsOpt.Rows.Item(17).Delete()
